$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("B","C","D","E","F","I","J","K","L","M","N")

$newValues = @{
    2 = @(1.02,1.041500726959166,1.042809449445845,1.049424201371748,1.059587858845765,1.041764316285448,1.046581597868908,1.045585088494984,1.052181298452591,1.062316954605946,1.019387048081298)
    3 = @(1.02,1.042522565608063,1.043571594603126,1.050405863304858,1.060819475880776,1.042060898866488,1.047249130345354,1.046158268481656,1.052974766166621,1.063361775068393,1.019612474797878)
    4 = @(1.02,1.043183739775688,1.044064642681008,1.051041449461118,1.061617288383875,1.042251419096339,1.04768045206977,1.046528382357188,1.053487961682206,1.064038150279682,1.019758039069796)
    5 = @(1.02,1.04346169180272,1.044271892983075,1.05130874192134,1.061952898269934,1.042331181209005,1.047861631959427,1.046683793137255,1.053703653637468,1.064322571943794,1.019819161873376)
    6 = @(1.02,1.043508360846009,1.044306689614007,1.05135362689341,1.062009260940857,1.042344554110281,1.047892044189715,1.046709876428451,1.053739866011954,1.064370331908627,1.019829420410342)
    7 = @(1.02,1.043187453806125,1.044067412076947,1.051045020674631,1.061621771995994,1.042252486188669,1.04768287358517,1.046530459691303,1.053490843986632,1.064041950446784,1.019758856080416)
    8 = @(1.02,1.041846067093779,1.043067042562321,1.049755878737874,1.060003908810836,1.041864835157104,1.046807321115418,1.045778957004006,1.05244950238254,1.062669993196586,1.019463294517077)
    9 = @(1.02,1.039482185499676,1.041303440422732,1.047487204308245,1.057159707052039,1.041171117500864,1.04525977272911,1.044448814202829,1.050612759133474,1.060254769746391,1.018940171405395)
    10 = @(1.02,1.037906129101758,1.040127184413451,1.045976751649598,1.055268042083044,1.04070150286699,1.044224913637805,1.0435581015646,1.049387079120883,1.058646181842939,1.018589880023473)
    11 = @(1.02,1.037223644064274,1.039617735110478,1.045323183420831,1.054449984627616,1.040496463355857,1.043776059786233,1.043171478729596,1.048856065035866,1.057950012093535,1.018437835488799)
    12 = @(1.02,1.036970132009385,1.039428484944344,1.045080489314586,1.054146278206922,1.040420048202943,1.043609222342134,1.043027728943868,1.048658779597467,1.057691477229178,1.018381304464833)
    13 = @(1.02,1.037024511516129,1.039469080555746,1.045132544846597,1.054211417190874,1.040436451021721,1.04364501466004,1.043058570132555,1.04870109994876,1.057746931417462,1.01839343304504)
    14 = @(1.02,1.037202688818107,1.039602091997718,1.045303120817767,1.05442487697019,1.040490152043851,1.043762271259066,1.043159599200012,1.048839758241161,1.057928640430464,1.018433163736481)
    15 = @(1.02,1.037312468832032,1.039684042362864,1.045408227608877,1.054556417310732,1.040523205313626,1.043834501927654,1.043221827882727,1.048925184468436,1.058040604437485,1.018457635852554)
    16 = @(1.02,1.037951422463899,1.040160992336921,1.046020136671898,1.055322355847913,1.040715074992636,1.044254686723782,1.043583740689376,1.049422314729898,1.05869239187304,1.018599963013108)
    17 = @(1.02,1.038352209623823,1.040460137905489,1.046404096260734,1.055803088464282,1.040834976467828,1.044518055933962,1.043810507892122,1.049734074611282,1.059101336841657,1.018689143135163)
    18 = @(1.02,1.03858597794565,1.040634612588725,1.046628098525516,1.056083592626136,1.040904749590036,1.044671602081576,1.043942686815593,1.049915891235006,1.059339902230836,1.018741125067125)
    19 = @(1.02,1.038665686227336,1.040694101921395,1.046704485165533,1.056179254439364,1.040928512708176,1.044723945026163,1.043987741049951,1.049977881320261,1.059421252763391,1.018758843568467)
    20 = @(1.02,1.038309209375875,1.040428043641275,1.046362896362488,1.055751499965565,1.040822129073015,1.044489806425378,1.04378618728228,1.049700628586175,1.059057457329666,1.018679578605212)
    21 = @(1.02,1.037150220232307,1.039562923937124,1.045252888507254,1.054362014069005,1.040474345452199,1.043727745236394,1.043129852544838,1.048798928022447,1.057875130179489,1.018421465553354)
    22 = @(1.02,1.036421478665244,1.039018884838717,1.044555388846577,1.053489294547307,1.040254208524468,1.043247952535431,1.042716373993114,1.048231743302129,1.057132063949666,1.018258861957919)
    23 = @(1.02,1.03680780219073,1.039307299976838,1.044925108110823,1.053951854055772,1.040371046746369,1.043502361772396,1.042935643991589,1.048532442368241,1.057525948194645,1.018345091285196)
    24 = @(1.02,1.038328639349296,1.040442545694413,1.046381512683486,1.055774810275325,1.040827934762471,1.044502571386706,1.043797176998701,1.049715741490518,1.059077284484904,1.018683900513718)
    25 = @(1.02,1.040093328698449,1.041759468073672,1.048073359629875,1.057894212019885,1.041351718592091,1.045660408712732,1.044793385475079,1.051087810292456,1.060878886839183,1.019075683732071)
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $col = $columns[$i]
        $ws.Range("$col$row").Value = $vals[$i]
    }
}
